# Applies the LOQ4241 syllabus content update:
#  - "Objetivos:" row gets real PT objective text (row 10 keeps the same
#    shared-string slot it always pointed at, so this is really just a
#    content fill-in for an already-referenced string).
#  - Rows 13-24 are re-shuffled: several rows lose their "A" label cell,
#    some gain new B/C content cells, the "Programa resumido:",
#    "Programa:", "Método:"/"Critério:", "Norma de recuperação:" and
#    "Bibliografia:" sections get real content, and a brand-new row 24
#    (Requisitos detail) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Long text blocks (kept as here-strings so accented/UTF-8 text and the
# embedded newlines come through exactly as in the source document).
# ---------------------------------------------------------------------

$objetivosTxt = "Apresentar conceitos, ferramentas e métodos para o auxílio à tomada de decisão."

$programaResumidoTxt = "Teoria da Decisão; Estruturação, Decisão sem Risco e sem Incerteza; Decisão com Múltiplos Cenários ou Múltiplos Critérios; Decisão com Incerteza; Sistemas de Auxílio à Decisão e Sistemas Especialistas."

$programaTxt = @"
1.Teoria da Decisão
2.Estruturação, Decisão sem Risco e sem Incerteza
3.Decisão com Múltiplos Cenários ou Múltiplos Critérios
4.Decisão com Incerteza
5.Sistemas de Auxílio à Decisão e Sistemas Especialistas.
"@

$metodoTxt = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."

$criterioTxt = "A Nota Final do aluno será determinada segundo a seguinte equação: Nota Final = (Prova- Bimestral-1*0,4) + (Prova-Bimestral-2*0,4) + (Trabalho*0,2)"

$normaRecuperacaoTxt = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."

$bibliografiaTxt = @"
1.ENSSLIN, L.; MONTIBELLER NETO, G.; NORONHA, S. M. (2001), Apoio à Decisão: metodologias para estruturação de problemas e avaliação multicritério de alternativas, Florianópolis: Insular
2.GOMES, L. F. A. M.; GOMES, C. F. S.; ALMEIDA, A. T. (2002), Tomada de Decisão Gerencial: enfoque multicritério, São Paulo: Atlas
3.LAWRENCE, J. A. JR.; PASTERNACK, B. A. (2002), Applied Management Science: modeling, spreadsheet analysis and communication for decision making, 2nd edition, New York (USA): Wiley
4.SHIMIZU, T. (2001), Decisão nas Organizações: introdução aos problemas de decisão encontrados nas organizações e nos sistemas de apoio à decisão, São Paulo: Atlas
"@

$requisitosDetailTxt = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"

# ---------------------------------------------------------------------
# "Objetivos:" (row 10) finally gets real content instead of reusing the
# professor-name string.
# ---------------------------------------------------------------------
$ws.Range("B10").Value = $objetivosTxt
$ws.Range("C10").Value = $objetivosTxt

# ---------------------------------------------------------------------
# Row 13: "Programa resumido:" label moves away (goes to row 14); row 13
# keeps only the professor-name value in B/C.
# ---------------------------------------------------------------------
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5840917 - Fabricio Maciel Gomes"
$ws.Range("C13").Value = "5840917 - Fabricio Maciel Gomes"
$ws.Rows.Item(13).AutoFit()

# ---------------------------------------------------------------------
# Row 14: now holds "Programa resumido:" + its real content.
# (B14 is a brand-new cell; column B's style ("2": top-aligned + wrap,
#  non-bold) is picked up by copy/paste-format from an existing column-B
#  data cell rather than toggled property-by-property, which avoids
#  leaving an orphaned intermediate style behind.)
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = $programaResumidoTxt
$ws.Range("C14").Value = $programaResumidoTxt
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 60

# ---------------------------------------------------------------------
# Row 15: becomes "Short syllabus:" label only (loses the stray
# B15/C15 date values it used to borrow).
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Rows.Item(15).RowHeight = 60

# ---------------------------------------------------------------------
# Row 16: "Programa:" gets its real (numbered) content.
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = $programaTxt
$ws.Range("C16").Value = $programaTxt
$ws.Range("B10").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Rows.Item(16).RowHeight = 120

# ---------------------------------------------------------------------
# Row 17: "Syllabus:" label only.
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Syllabus:"
$ws.Rows.Item(17).RowHeight = 120

# ---------------------------------------------------------------------
# Row 18: "Avaliação:" label only now (loses the borrowed
# professor-name B18/C18 values).
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).AutoFit()

# ---------------------------------------------------------------------
# Row 19: "Método:" gets its real content (moved up from row 20 before).
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = $metodoTxt
$ws.Range("C19").Value = $metodoTxt
$ws.Rows.Item(19).RowHeight = 60

# ---------------------------------------------------------------------
# Row 20: "Critério:" gets the grading-formula content.
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = $criterioTxt
$ws.Range("C20").Value = $criterioTxt
$ws.Rows.Item(20).RowHeight = 60

# ---------------------------------------------------------------------
# Row 21: "Norma de recuperação:" gets the make-up exam norm text.
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = $normaRecuperacaoTxt
$ws.Range("C21").Value = $normaRecuperacaoTxt
$ws.Rows.Item(21).RowHeight = 60

# ---------------------------------------------------------------------
# Row 22: "Bibliografia:" gets the real reading list.
# ---------------------------------------------------------------------
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = $bibliografiaTxt
$ws.Range("C22").Value = $bibliografiaTxt
$ws.Range("B10").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Rows.Item(22).RowHeight = 120

# ---------------------------------------------------------------------
# Row 23: "Requisitos:" label only now (its detail text moves to the new
# row 24).
# ---------------------------------------------------------------------
$ws.Range("A23").Value = "Requisitos:"
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()
$ws.Rows.Item(23).AutoFit()

# ---------------------------------------------------------------------
# Row 24 (brand new): requisite detail text in B/C.
# ---------------------------------------------------------------------
$ws.Range("B24").Value = $requisitosDetailTxt
$ws.Range("C24").Value = $requisitosDetailTxt
$ws.Range("B10").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Rows.Item(24).RowHeight = 30
